$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H70").Value = 967.7273
$ws.Range("I70").Value = 938.4
$ws.Range("J70").Value = 992.1667
$ws.Range("K70").Value = 2815.2
$ws.Range("L70").Value = 2976.5001
$ws.Range("M70").Value = -2545.2
$ws.Range("N70").Value = -3516.5001

$ws.Range("H73").Value = 967.7273
$ws.Range("I73").Value = 938.4
$ws.Range("J73").Value = 992.1667
$ws.Range("K73").Value = 2815.2
$ws.Range("L73").Value = 2976.5001
$ws.Range("M73").Value = -1879.2
$ws.Range("N73").Value = -4848.5001

$ws.Range("H76").Value = 3190
$ws.Range("I76").Value = 3250
$ws.Range("J76").Value = 2950
$ws.Range("K76").Value = 3250
$ws.Range("L76").Value = 2950
$ws.Range("M76").Value = -2935
$ws.Range("N76").Value = -3580

$ws.Range("H79").Value = 3190
$ws.Range("I79").Value = 3250
$ws.Range("J79").Value = 2950
$ws.Range("K79").Value = 3250
$ws.Range("L79").Value = 2950
$ws.Range("M79").Value = -2158
$ws.Range("N79").Value = -5134

$ws.Range("H92").Value = 66670436
$ws.Range("I92").Value = 83337040
$ws.Range("J92").Value = 4000
$ws.Range("K92").Value = 83337040
$ws.Range("L92").Value = 4000
$ws.Range("M92").Value = -83335792
$ws.Range("N92").Value = -6496

$ws.Range("H98").Value = 554.1111
$ws.Range("I98").Value = 248.375
$ws.Range("K98").Value = 248.375
$ws.Range("M98").Value = 1249.625

$ws.Range("H112").Value = 4835.4614
$ws.Range("J112").Value = 5219.4585
$ws.Range("L112").Value = 15658.3755
$ws.Range("N112").Value = -17874.3755

$ws.Range("H118").Value = 2390
$ws.Range("I118").Value = 2455
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 7365
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -5708
$ws.Range("N118").Value = -9314

$ws.Range("H122").Value = 554.1111
$ws.Range("I122").Value = 248.375
$ws.Range("K122").Value = 745.125
$ws.Range("M122").Value = 1704.875

$ws.Range("H137").Value = 3307.5518
$ws.Range("I137").Value = 2800.6956
$ws.Range("J137").Value = 5250.5
$ws.Range("K137").Value = 8402.086800000001
$ws.Range("L137").Value = 15751.5
$ws.Range("M137").Value = -5852.086800000001
$ws.Range("N137").Value = -20851.5

$ws.Range("H138").Value = 2220.15
$ws.Range("I138").Value = 1607.5
$ws.Range("J138").Value = 2756.2188
$ws.Range("K138").Value = 4822.5
$ws.Range("L138").Value = 8268.6564
$ws.Range("M138").Value = 317.5
$ws.Range("N138").Value = -18548.6564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10102.6
$ws.Range("I28").Value = 10102.6
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 10102.6
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -9910.6
$ws.Range("N28").ClearContents()

$ws.Range("H32").Value = 361695.9
$ws.Range("I32").Value = 386830.97
$ws.Range("K32").Value = 386830.97
$ws.Range("M32").Value = -386543.97

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 6214.5
$ws.Range("I63").Value = 6287.143
$ws.Range("J63").Value = 6141.857
$ws.Range("K63").Value = 6287.143
$ws.Range("L63").Value = 6141.857
$ws.Range("M63").Value = -5601.143
$ws.Range("N63").Value = -7513.857

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 6214.5
$ws.Range("I66").Value = 6287.143
$ws.Range("J66").Value = 6141.857
$ws.Range("K66").Value = 31435.715
$ws.Range("L66").Value = 30709.285
$ws.Range("M66").Value = -28003.715
$ws.Range("N66").Value = -37573.285

$ws.Range("H97").Value = 594.9
$ws.Range("I97").Value = 594.9
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 594.9
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -98.89999999999998
$ws.Range("N97").ClearContents()

$ws.Range("H99").Value = 10102.6
$ws.Range("I99").Value = 10102.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 10102.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -7107.6
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1257
$ws.Range("I94").Value = 1102.25
$ws.Range("J94").Value = 1463.3334
$ws.Range("K94").Value = 1102.25
$ws.Range("L94").Value = 1463.3334
$ws.Range("M94").Value = -651.25
$ws.Range("N94").Value = -2365.3334

$ws.Range("H105").Value = 2724.32
$ws.Range("I105").Value = 2700.348
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2700.348
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -953.348
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 27562.562
$ws.Range("J4").Value = 27733.4
$ws.Range("L4").Value = 27733.4
$ws.Range("N4").Value = -27957.4

$ws.Range("H22").Value = 543
$ws.Range("I22").Value = 481.81818
$ws.Range("J22").Value = 767.3333
$ws.Range("K22").Value = 481.81818
$ws.Range("L22").Value = 767.3333
$ws.Range("M22").Value = -131.81818
$ws.Range("N22").Value = -1467.3333

$ws.Range("H31").Value = 5453.9165
$ws.Range("I31").Value = 1351.8
$ws.Range("J31").Value = 8384
$ws.Range("K31").Value = 1351.8
$ws.Range("L31").Value = 8384
$ws.Range("M31").Value = -1056.8
$ws.Range("N31").Value = -8974

$ws.Range("H34").Value = 5453.9165
$ws.Range("I34").Value = 1351.8
$ws.Range("J34").Value = 8384
$ws.Range("K34").Value = 1351.8
$ws.Range("L34").Value = 8384
$ws.Range("M34").Value = -1149.8
$ws.Range("N34").Value = -8788

$ws.Range("H105").Value = 2499.75
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2999.5
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2999.5
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6493.5

$ws.Range("H116").Value = 63000
$ws.Range("J116").Value = 63000
$ws.Range("L116").Value = 63000
$ws.Range("N116").Value = -72178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 155.57143
$ws.Range("I2").Value = 89.666664
$ws.Range("K2").Value = 537.999984
$ws.Range("M2").Value = -424.999984

$ws.Range("H4").Value = 18463730
$ws.Range("J4").Value = 20002400
$ws.Range("L4").Value = 60007200
$ws.Range("N4").Value = -60007424

$ws.Range("H92").Value = 595
$ws.Range("I92").Value = 590
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 1770
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = -522
$ws.Range("N92").Value = -4296

$ws.Range("H97").Value = 100004
$ws.Range("J97").Value = 100004
$ws.Range("L97").Value = 300012
$ws.Range("N97").Value = -301004

$ws.Range("H98").Value = 301
$ws.Range("I98").Value = 150
$ws.Range("J98").Value = 351.33334
$ws.Range("K98").Value = 450
$ws.Range("L98").Value = 1054.00002
$ws.Range("M98").Value = 1048
$ws.Range("N98").Value = -4050.00002

$ws.Range("H122").Value = 10195
$ws.Range("J122").Value = 15570.714
$ws.Range("L122").Value = 140136.426
$ws.Range("N122").Value = -145036.426

$ws.Range("H136").Value = 1481.5
$ws.Range("I136").Value = 977.8
$ws.Range("K136").Value = 2933.4
$ws.Range("M136").Value = 2166.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 72773000
$ws.Range("I80").Value = 127251256
$ws.Range("J80").Value = 135333
$ws.Range("K80").Value = 127251256
$ws.Range("L80").Value = 135333
$ws.Range("M80").Value = -127250258
$ws.Range("N80").Value = -137329

$ws.Range("H83").Value = 72773000
$ws.Range("I83").Value = 127251256
$ws.Range("J83").Value = 135333
$ws.Range("K83").Value = 636256280
$ws.Range("L83").Value = 676665
$ws.Range("M83").Value = -636251288
$ws.Range("N83").Value = -686649

$ws.Range("H97").Value = 1912.2727
$ws.Range("I97").Value = 1732.3077
$ws.Range("J97").Value = 2172.2222
$ws.Range("K97").Value = 1732.3077
$ws.Range("L97").Value = 2172.2222
$ws.Range("M97").Value = -1236.3077
$ws.Range("N97").Value = -3164.2222

$ws.Range("H113").Value = 168070.5
$ws.Range("I113").Value = 250605.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 250605.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -248435.75
$ws.Range("N113").Value = -7340

$ws.Range("H114").Value = 35637
$ws.Range("J114").Value = 35637
$ws.Range("L114").Value = 35637
$ws.Range("N114").Value = -44315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3380001.2
$ws.Range("J2").Value = 3380001.2
$ws.Range("L2").Value = 3380001.2
$ws.Range("N2").Value = -3380225.2

$ws.Range("H18").Value = 57504.5
$ws.Range("J18").Value = 57504.5
$ws.Range("L18").Value = 57504.5
$ws.Range("N18").Value = -57848.5

$ws.Range("H61").Value = 3181.6667
$ws.Range("I61").Value = 1782.5
$ws.Range("K61").Value = 1782.5
$ws.Range("M61").Value = -1580.5

$ws.Range("H113").Value = 3181.6667
$ws.Range("I113").Value = 1782.5
$ws.Range("K113").Value = 1782.5
$ws.Range("M113").Value = 387.5

$ws.Range("H132").Value = 2849.9583
$ws.Range("I132").Value = 1996.4166
$ws.Range("J132").Value = 3703.5
$ws.Range("K132").Value = 5989.2498
$ws.Range("L132").Value = 11110.5
$ws.Range("M132").Value = -3459.2498
$ws.Range("N132").Value = -16170.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3358667.8
$ws.Range("J2").Value = 38001.5
$ws.Range("L2").Value = 38001.5
$ws.Range("N2").Value = -38225.5

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
